$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 59: Melón / Tuna / Extra
$ws.Cells.Item(59, 1).Value = 8
$ws.Cells.Item(59, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(59, 3).Value = "Coquimbo"
$ws.Cells.Item(59, 4).Value = 44568
$ws.Cells.Item(59, 4).NumberFormat = $ws.Cells.Item(58, 4).NumberFormat
$ws.Cells.Item(59, 5).Value = 4
$ws.Cells.Item(59, 6).Value = 100112027
$ws.Cells.Item(59, 7).Value = "Melón"
$ws.Cells.Item(59, 8).Value = "Tuna"
$ws.Cells.Item(59, 9).Value = "Extra"
$ws.Cells.Item(59, 10).Value = 6000
$ws.Cells.Item(59, 11).Value = 1100
$ws.Cells.Item(59, 12).Value = 1200
$ws.Cells.Item(59, 13).Value = 1150
$ws.Cells.Item(59, 14).Value = "$/unidad"
$ws.Cells.Item(59, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(59, 16).Value = 1150
$ws.Cells.Item(59, 17).Value = 1
$ws.Cells.Item(59, 18).Value = "Hortaliza"

# Row 60: Melón / Tuna / Primera
$ws.Cells.Item(60, 1).Value = 8
$ws.Cells.Item(60, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(60, 3).Value = "Coquimbo"
$ws.Cells.Item(60, 4).Value = 44568
$ws.Cells.Item(60, 4).NumberFormat = $ws.Cells.Item(58, 4).NumberFormat
$ws.Cells.Item(60, 5).Value = 4
$ws.Cells.Item(60, 6).Value = 100112027
$ws.Cells.Item(60, 7).Value = "Melón"
$ws.Cells.Item(60, 8).Value = "Tuna"
$ws.Cells.Item(60, 9).Value = "Primera"
$ws.Cells.Item(60, 10).Value = 5000
$ws.Cells.Item(60, 11).Value = 950
$ws.Cells.Item(60, 12).Value = 1000
$ws.Cells.Item(60, 13).Value = 975
$ws.Cells.Item(60, 14).Value = "$/unidad"
$ws.Cells.Item(60, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(60, 16).Value = 975
$ws.Cells.Item(60, 17).Value = 1
$ws.Cells.Item(60, 18).Value = "Hortaliza"
